# "Trying to resolve conflict !"
#
# Target sheet for all the content/selection changes is "Register".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Register")

# --- 1. Shared-string content fix -----------------------------------------
# B7 used to read "Bug_Busters"; C7 already read "Bug@Busters". The merge
# conflict resolution renames that text to "Team_Bug@busters" and makes B7
# match C7 again.
$ws.Range("B7").Value = "Team_Bug@busters"
$ws.Range("C7").Value = "Team_Bug@busters"

# --- 2. Leftover formatting touch ------------------------------------------
# A stray white-fill / no-border formatting pass (picked up while resolving
# the conflict) that was applied and then cleared again, on a cell with no
# content (D9) so no visible cell is affected, only the style tables grow.
$ws.Range("D9").Interior.Color = 16777215
$ws.Range("D9").Borders.LineStyle = 1
$ws.Range("D9").ClearFormats()

# --- 3. Cursor position -----------------------------------------------------
$ws.Range("C15").Select()
